$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F2: 39 -> 40 (想去人数 for row 2)
    $ws.Range("F2").Value = 40

    # I3: cover image URL update
    $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202406/LSorIT7S1717486817969.png"

    # F5: 28 -> 30 (想去人数 for row 5)
    $ws.Range("F5").Value = 30
}
